# Auto-generated edit script
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H4").Value = 385.06668
$ws_ALC.Range("I4").Value = 248.8
$ws_ALC.Range("K4").Value = 248.8
$ws_ALC.Range("M4").Value = -134.8

$ws_ALC.Range("H40").Value = 2916.6667
$ws_ALC.Range("I40").Value = 2500
$ws_ALC.Range("K40").Value = 2500
$ws_ALC.Range("M40").Value = -2325

$ws_ALC.Range("H55").Value = 309.5
$ws_ALC.Range("J55").Value = 385.33334
$ws_ALC.Range("L55").Value = 385.33334
$ws_ALC.Range("N55").Value = -813.33334

$ws_ALC.Range("H86").Value = 3218.647
$ws_ALC.Range("I86").Value = 3655.25
$ws_ALC.Range("K86").Value = 3655.25
$ws_ALC.Range("M86").Value = -2532.25

$ws_ALC.Range("H89").Value = 3218.647
$ws_ALC.Range("I89").Value = 3655.25
$ws_ALC.Range("K89").Value = 18276.25
$ws_ALC.Range("M89").Value = -12660.25

$ws_ALC.Range("H98").Value = 904.1667
$ws_ALC.Range("J98").Value = 734.1429000000001
$ws_ALC.Range("L98").Value = 734.1429000000001
$ws_ALC.Range("N98").Value = -3730.1429

$ws_ALC.Range("H107").Value = 3271
$ws_ALC.Range("I107").Value = 4380
$ws_ALC.Range("K107").Value = 4380
$ws_ALC.Range("M107").Value = -2460

$ws_ALC.Range("H122").Value = 904.1667
$ws_ALC.Range("J122").Value = 734.1429000000001
$ws_ALC.Range("L122").Value = 2202.4287
$ws_ALC.Range("N122").Value = -7102.4287

$ws_ALC.Range("H132").Value = 3153.7742
$ws_ALC.Range("I132").Value = 2123.55
$ws_ALC.Range("K132").Value = 6370.650000000001
$ws_ALC.Range("M132").Value = -3840.650000000001

$ws_ALC.Range("H137").Value = 2007.3846
$ws_ALC.Range("I137").Value = 1490.0952
$ws_ALC.Range("K137").Value = 4470.2856
$ws_ALC.Range("M137").Value = -1920.2856

$ws_ALC.Range("H138").Value = 3478.625
$ws_ALC.Range("J138").Value = 3728.4285
$ws_ALC.Range("L138").Value = 11185.2855
$ws_ALC.Range("N138").Value = -21465.2855

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 1800.2858
$ws_ARM.Range("I2").Value = 431.83334
$ws_ARM.Range("J2").Value = 10011
$ws_ARM.Range("K2").Value = 431.83334
$ws_ARM.Range("L2").Value = 10011
$ws_ARM.Range("M2").Value = -318.83334
$ws_ARM.Range("N2").Value = -10237

$ws_ARM.Range("H32").Value = 6857.263
$ws_ARM.Range("I32").Value = 5018.0625
$ws_ARM.Range("K32").Value = 5018.0625
$ws_ARM.Range("M32").Value = -4731.0625

$ws_ARM.Range("H116").Value = 1800.2858
$ws_ARM.Range("I116").Value = 431.83334
$ws_ARM.Range("J116").Value = 10011
$ws_ARM.Range("K116").Value = 431.83334
$ws_ARM.Range("L116").Value = 10011
$ws_ARM.Range("M116").Value = 1862.16666
$ws_ARM.Range("N116").Value = -14599

$ws_ARM.Range("H132").Value = 1307
$ws_ARM.Range("I132").Value = 1237.7
$ws_ARM.Range("K132").Value = 3713.1
$ws_ARM.Range("M132").Value = -1183.1

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 1800.2858
$ws_BSM.Range("I3").Value = 431.83334
$ws_BSM.Range("J3").Value = 10011
$ws_BSM.Range("K3").Value = 431.83334
$ws_BSM.Range("L3").Value = 10011
$ws_BSM.Range("M3").Value = -317.83334
$ws_BSM.Range("N3").Value = -10239

$ws_BSM.Range("H105").Value = 2231
$ws_BSM.Range("J105").Value = 5750
$ws_BSM.Range("L105").Value = 5750
$ws_BSM.Range("N105").Value = -9244

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3415
$ws_CRP.Range("I31").Value = 3415
$ws_CRP.Range("K31").Value = 3415
$ws_CRP.Range("M31").Value = -3120

$ws_CRP.Range("H34").Value = 3415
$ws_CRP.Range("I34").Value = 3415
$ws_CRP.Range("K34").Value = 3415
$ws_CRP.Range("M34").Value = -3213

$ws_CRP.Range("H58").Value = 1236.1666
$ws_CRP.Range("I58").Value = 937.4286
$ws_CRP.Range("K58").Value = 937.4286
$ws_CRP.Range("M58").Value = -734.4286

$ws_CRP.Range("H122").Value = 1796.7059
$ws_CRP.Range("I122").Value = 1004.125
$ws_CRP.Range("K122").Value = 3012.375
$ws_CRP.Range("M122").Value = -562.375

$ws_CRP.Range("H136").Value = 1236.1666
$ws_CRP.Range("I136").Value = 937.4286
$ws_CRP.Range("K136").Value = 2812.2858
$ws_CRP.Range("M136").Value = -262.2857999999997

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H3").Value = 520
$ws_CUL.Range("I3").Value = 200
$ws_CUL.Range("K3").Value = 600
$ws_CUL.Range("M3").Value = -488

$ws_CUL.Range("H34").Value = 837.5
$ws_CUL.Range("J34").Value = 994.3333
$ws_CUL.Range("L34").Value = 2982.9999
$ws_CUL.Range("N34").Value = -3150.9999

$ws_CUL.Range("H55").Value = 1999.5
$ws_CUL.Range("I55").Value = 0
$ws_CUL.Range("K55").Value = 0
$ws_CUL.Range("M55").ClearContents()

$ws_CUL.Range("H94").Value = 0
$ws_CUL.Range("I94").Value = 0
$ws_CUL.Range("K94").Value = 0
$ws_CUL.Range("M94").ClearContents()

$ws_CUL.Range("H121").Value = 1502.8334
$ws_CUL.Range("J121").Value = 2011.875
$ws_CUL.Range("L121").Value = 6035.625
$ws_CUL.Range("N121").Value = -8655.625

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H33").Value = 0
$ws_GSM.Range("I33").Value = 0
$ws_GSM.Range("K33").Value = 0
$ws_GSM.Range("M33").ClearContents()

$ws_GSM.Range("H80").Value = 3150
$ws_GSM.Range("I80").Value = 0
$ws_GSM.Range("J80").Value = 3150
$ws_GSM.Range("K80").Value = 0
$ws_GSM.Range("L80").Value = 3150
$ws_GSM.Range("M80").ClearContents()
$ws_GSM.Range("N80").Value = -5146

$ws_GSM.Range("H83").Value = 3150
$ws_GSM.Range("I83").Value = 0
$ws_GSM.Range("J83").Value = 3150
$ws_GSM.Range("K83").Value = 0
$ws_GSM.Range("L83").Value = 15750
$ws_GSM.Range("M83").ClearContents()
$ws_GSM.Range("N83").Value = -25734

$ws_GSM.Range("H97").Value = 1073.7
$ws_GSM.Range("I97").Value = 693
$ws_GSM.Range("K97").Value = 693
$ws_GSM.Range("M97").Value = -197

$ws_GSM.Range("H113").Value = 2100
$ws_GSM.Range("J113").Value = 2000
$ws_GSM.Range("L113").Value = 2000
$ws_GSM.Range("N113").Value = -6340

$ws_GSM.Range("H132").Value = 1681.25
$ws_GSM.Range("I132").Value = 1761.3636
$ws_GSM.Range("J132").Value = 800
$ws_GSM.Range("K132").Value = 5284.0908
$ws_GSM.Range("L132").Value = 2400
$ws_GSM.Range("M132").Value = -2754.0908
$ws_GSM.Range("N132").Value = -7460

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 2219
$ws_LTW.Range("I7").Value = 1711.8572
$ws_LTW.Range("K7").Value = 1711.8572
$ws_LTW.Range("M7").Value = -1599.8572

$ws_LTW.Range("H22").Value = 5150
$ws_LTW.Range("I22").Value = 2300
$ws_LTW.Range("K22").Value = 2300
$ws_LTW.Range("M22").Value = -2005

$ws_LTW.Range("H27").Value = 5150
$ws_LTW.Range("I27").Value = 2300
$ws_LTW.Range("K27").Value = 2300
$ws_LTW.Range("M27").Value = -2193

$ws_LTW.Range("H126").Value = 2219
$ws_LTW.Range("I126").Value = 1711.8572
$ws_LTW.Range("K126").Value = 5135.571599999999
$ws_LTW.Range("M126").Value = -2665.571599999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H107").Value = 300.63635
$ws_WVR.Range("I107").Value = 354.66666
$ws_WVR.Range("K107").Value = 1063.99998
$ws_WVR.Range("M107").Value = 856.0000199999999

$ws_WVR.Range("H126").Value = 2600.4614
$ws_WVR.Range("I126").Value = 2704.6667
$ws_WVR.Range("K126").Value = 8114.000100000001
$ws_WVR.Range("M126").Value = -5644.000100000001
